$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" / Changed date column) for data rows 2-28:
# each value moves from 45530 to 45531 (one day later), keeping formatting intact.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
